$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (x_nrSteps) was stored as negative zero (-0) for every data row;
# normalize it to a plain 0 for rows 2-31.
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
}

# A handful of rows had their y_corrSteps (E), y_nrSteps (G) and alienID (H)
# values recomputed; apply the corrected numbers.
$rowUpdates = @{
    4  = @(6, 3, 13)
    8  = @(6, 3, 13)
    16 = @(7, 3, 13)
    18 = @(6, 3, 13)
    23 = @(5, 3, 13)
    27 = @(7, 3, 13)
}

foreach ($r in $rowUpdates.Keys) {
    $vals = $rowUpdates[$r]
    $ws.Cells.Item($r, 5).Value = $vals[0]
    $ws.Cells.Item($r, 7).Value = $vals[1]
    $ws.Cells.Item($r, 8).Value = $vals[2]
}
